# Updates cryptos price list per latest data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.990.78'
$ws.Range('E2').Value = '  -0.18%  '
$ws.Range('D3').Value = '2.040.80'
$ws.Range('E3').Value = '  -0.31%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '247.91'
$ws.Range('E5').Value = '  -0.19%  '
$ws.Range('E6').Value = '  +0.68%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '56.24'
$ws.Range('E8').Value = '  +1.40%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.380'
$ws.Range('E9').Value = '  +0.29%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0779'
$ws.Range('E10').Value = '  +2.94%  '
$ws.Range('E11').Value = '  +1.67%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.74'
$ws.Range('E12').Value = '  +4.47%  '
$ws.Range('D13').Value = '2.337.57'
$ws.Range('E13').Value = '  -0.50%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.60'
$ws.Range('E14').Value = '  +6.63%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.788'
$ws.Range('E15').Value = '  -4.42%  '
$ws.Range('D16').Value = '2.039.10'
$ws.Range('E16').Value = '  -1.04%  '
$ws.Range('D17').Value = '36.937.73'
$ws.Range('E17').Value = '  -0.19%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '16.53'
$ws.Range('E18').Value = '  +14.19%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '73.71'
$ws.Range('E19').Value = '  +1.97%  '
$ws.Range('D20').Value = '0.0₃0893'
$ws.Range('E20').Value = '  +4.21%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.30'
$ws.Range('E21').Value = '  +1.37%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '235.29'
$ws.Range('E22').Value = '  -1.09%  '
$ws.Range('E23').Value = '  +0.08%  '
$ws.Range('E24').Value = '  -2.10%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.19'
$ws.Range('E25').Value = '  +9.50%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '167.52'
$ws.Range('E26').Value = '  -1.46%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.08'
$ws.Range('E27').Value = '  -0.26%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.66'
$ws.Range('E28').Value = '  -3.38%  '
$ws.Range('E29').Value = '  +0.99%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.11'
$ws.Range('E30').Value = '  +2.54%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.67'
$ws.Range('E31').Value = '  +2.52%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0608'
$ws.Range('E32').Value = '  -3.45%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.41'
$ws.Range('E33').Value = '  +1.27%  '
$ws.Range('E34').Value = '  -0.21%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0867'
$ws.Range('E35').Value = '  -0.10%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.20'
$ws.Range('E36').Value = '  -2.63%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.75'
$ws.Range('E37').Value = '  -3.15%  '
$ws.Range('B38').Value = 'Cronos'
$ws.Range('C38').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.107'
$ws.Range('E38').Value = '  +1.06%  '
$ws.Range('B39').Value = 'HuobiToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.24'
$ws.Range('E39').Value = '  +15.35%  '
$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.33'
$ws.Range('E40').Value = '  -0.74%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0220'
$ws.Range('E41').Value = '  -2.23%  '
$ws.Range('B42').Value = 'InjectiveProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '17.28'
$ws.Range('E42').Value = '  -4.70%  '
$ws.Range('B43').Value = 'THORChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.85'
$ws.Range('E43').Value = '  +22.10%  '
$ws.Range('B44').Value = 'ARBITRUM'
$ws.Range('C44').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.11'
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '95.19'
$ws.Range('E45').Value = '  -2.17%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.42'
$ws.Range('E46').Value = '  +1.71%  '
$ws.Range('D47').Value = '1.273.29'
$ws.Range('E47').Value = '  -2.34%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.84'
$ws.Range('E48').Value = '  -2.29%  '
$ws.Range('B49').Value = 'RocketPoolETH'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D49').Value = '2.224.09'
$ws.Range('E49').Value = '  -0.53%  '
$ws.Range('B50').Value = 'FraxShare'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.66'
$ws.Range('E50').Value = '  -4.22%  '
$ws.Range('B51').Value = 'FTXToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.39'
$ws.Range('E51').Value = '  -21.34%  '
